# Update countries & provincias Spain
#
# Refreshes the "Pais" sheet with a later data pull (17:35 -> 18:05) and
# re-applies the table's descending sort on "Casos totales" where the new
# counts changed the ranking:
#   - Polonia overtakes Banglades (rows 33/34 swap)
#   - Yemen overtakes Libia / Polinesia Francesa / Malaui (rows 171-174 shift)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 12 de Mayo de 2020 a las 18:05"

# row => column letter => new value (string for col A, number otherwise)
$updates = [ordered]@{
    4   = @{ B = 1390816; C = 4982; D = 262880; E = 1045814;            G = 327; H = 82122 }
    11  = @{ B = 170582;  C = 1439;                E = 91426;            G = 147; H = 11772 }
    30  = @{                        D = 3851;    E = 20799; F = 20 }
    33  = @{ A = "Polonia";             B = 16921; C = 595; D = 6131; E = 9951; F = 160; G = 28; H = 839 }
    34  = @{ A = "Banglades";           B = 16660; C = 969; D = 3147; E = 13263; F = 1;   G = 11; H = 250 }
    39  = @{                                         E = 7091;                   G = 20; H = 1002 }
    45  = @{                        D = 3221;    E = 7277 }
    56  = @{                        D = 1862;    E = 4099 }
    65  = @{ B = 3894; C = 6; D = 3610; E = 182; F = 22; G = 1; H = 102 }
    121 = @{ B = 576;  C = 14;               E = 177 }
    171 = @{ A = "Yemen";               B = 65; C = 9; D = 1;  E = 54;  G = 1; H = 10 }
    172 = @{ A = "Libia";               B = 64;        D = 28; E = 33;  F = 0;        H = 3 }
    173 = @{ A = "Polinesia Francesa";  B = 60;        D = 56; E = 4;                 H = 0 }
    174 = @{ A = "Malaui";              B = 57;        D = 24; E = 30;  F = 1;        H = 3 }
}

$colIndex = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8 }

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item([int]$row, $colIndex[$col]).Value = $rowData[$col]
    }
}
